$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.490.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.139.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "613.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.389"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.92%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.132.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.15%  "

$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.665.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000243"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.713.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.114.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("E19").Value = "  -5.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +16.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000195"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "86.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.290.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.242"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.02%  "

$ws.Range("E31").Value = "  -1.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.86%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.69%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.846"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -15.31%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.152"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.45%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "490.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.441"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "

$ws.Range("E43").Value = "  -6.85%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "162.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.73%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.704"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0326"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.53%  "
